$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update B/C values that changed ---
# Row 3: Temperature
$ws.Range("B3").Value = 19

# Row 4: Hour
$ws.Range("B4").Value = "16:50:29"
$ws.Range("C4").Value = "16:50:35"

# Row 5: Date
$ws.Range("B5").Value = "29-03-23"
$ws.Range("C5").Value = "29-03-23"

# --- Clear out columns D:U for the data rows (keep styles, drop content) ---
$ws.Range("D2:U2").ClearContents()
$ws.Range("D3:U3").ClearContents()
$ws.Range("D4:U4").ClearContents()
$ws.Range("D5:U5").ClearContents()
$ws.Range("D8:U8").ClearContents()
$ws.Range("D9:U9").ClearContents()
$ws.Range("D10:U10").ClearContents()
